$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 58.25
$ws.Range("I5").Value = 63
$ws.Range("J5").Value = 44
$ws.Range("K5").Value = 63
$ws.Range("L5").Value = 44
$ws.Range("M5").Value = 52
$ws.Range("N5").Value = -274

# Row 9
$ws.Range("H9").Value = 298
$ws.Range("I9").Value = 227.2
$ws.Range("J9").Value = 439.6
$ws.Range("K9").Value = 227.2
$ws.Range("L9").Value = 439.6
$ws.Range("M9").Value = -58.19999999999999
$ws.Range("N9").Value = -777.6

# Row 18
$ws.Range("H18").Value = 4899.6284
$ws.Range("I18").Value = 3819.48
$ws.Range("J18").Value = 7600
$ws.Range("K18").Value = 3819.48
$ws.Range("L18").Value = 7600
$ws.Range("M18").Value = -3535.48
$ws.Range("N18").Value = -8168

# Row 29
$ws.Range("H29").Value = 187.5
$ws.Range("I29").Value = 187.5
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 562.5
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -281.5
$ws.Range("N29").ClearContents()

# Row 33
$ws.Range("H33").Value = 553.5
$ws.Range("I33").Value = 497.92856
$ws.Range("J33").Value = 683.1667
$ws.Range("K33").Value = 497.92856
$ws.Range("L33").Value = 683.1667
$ws.Range("M33").Value = -268.92856
$ws.Range("N33").Value = -1141.1667

# Row 34
$ws.Range("H34").Value = 2532.625
$ws.Range("I34").Value = 1876.8334
$ws.Range("J34").Value = 4500
$ws.Range("K34").Value = 1876.8334
$ws.Range("L34").Value = 4500
$ws.Range("M34").Value = -1673.8334
$ws.Range("N34").Value = -4906

# Row 36
$ws.Range("H36").Value = 2532.625
$ws.Range("I36").Value = 1876.8334
$ws.Range("J36").Value = 4500
$ws.Range("K36").Value = 1876.8334
$ws.Range("L36").Value = 4500
$ws.Range("M36").Value = -1161.8334
$ws.Range("N36").Value = -5930

# Row 38
$ws.Range("H38").Value = 293
$ws.Range("I38").Value = 293
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 879
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -507
$ws.Range("N38").ClearContents()

# Row 43
$ws.Range("H43").Value = 5499.8184
$ws.Range("I43").Value = 8000
$ws.Range("J43").Value = 4944.222
$ws.Range("K43").Value = 8000
$ws.Range("L43").Value = 4944.222
$ws.Range("M43").Value = -7931
$ws.Range("N43").Value = -5082.222

# Row 58
$ws.Range("H58").Value = 794.7
$ws.Range("I58").Value = 99
$ws.Range("J58").Value = 1092.8572
$ws.Range("K58").Value = 297
$ws.Range("L58").Value = 3278.5716
$ws.Range("M58").Value = -147
$ws.Range("N58").Value = -3578.5716

# Row 113
$ws.Range("H113").Value = 9683.464
$ws.Range("I113").Value = 4915.6665
$ws.Range("J113").Value = 11941.895
$ws.Range("K113").Value = 4915.6665
$ws.Range("L113").Value = 11941.895
$ws.Range("M113").Value = -1661.6665
$ws.Range("N113").Value = -18449.895

# Row 116
$ws.Range("H116").Value = 7617.3887
$ws.Range("I116").Value = 6108.5
$ws.Range("J116").Value = 8824.5
$ws.Range("K116").Value = 6108.5
$ws.Range("L116").Value = 8824.5
$ws.Range("M116").Value = -2666.5
$ws.Range("N116").Value = -15708.5

# Row 126
$ws.Range("H126").Value = 39642.855
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 39642.855
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 39642.855
$ws.Range("N126").Value = -49522.855

# Row 128
$ws.Range("H128").Value = 39769.23
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 39769.23
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 39769.23
$ws.Range("N128").Value = -49729.23

# Row 137
$ws.Range("H137").Value = 4579.857
$ws.Range("I137").Value = 5678.091
$ws.Range("J137").Value = 3371.8
$ws.Range("K137").Value = 17034.273
$ws.Range("L137").Value = 10115.4
$ws.Range("M137").Value = -14484.273
$ws.Range("N137").Value = -15215.4


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 2647.4092
$ws.Range("I45").Value = 2172.7693
$ws.Range("J45").Value = 3333
$ws.Range("K45").Value = 2172.7693
$ws.Range("L45").Value = 3333
$ws.Range("M45").Value = -1795.7693
$ws.Range("N45").Value = -4087

# Row 74
$ws.Range("H74").Value = 1251.6
$ws.Range("I74").Value = 1265.6
$ws.Range("J74").Value = 1209.6
$ws.Range("K74").Value = 1265.6
$ws.Range("L74").Value = 1209.6
$ws.Range("M74").Value = -391.5999999999999
$ws.Range("N74").Value = -2957.6

# Row 77
$ws.Range("H77").Value = 1251.6
$ws.Range("I77").Value = 1265.6
$ws.Range("J77").Value = 1209.6
$ws.Range("K77").Value = 6328
$ws.Range("L77").Value = 6048
$ws.Range("M77").Value = -1960
$ws.Range("N77").Value = -14784

# Row 110
$ws.Range("H110").Value = 3311.8147
$ws.Range("I110").Value = 2815.35
$ws.Range("J110").Value = 4730.2856
$ws.Range("K110").Value = 2815.35
$ws.Range("L110").Value = 4730.2856
$ws.Range("M110").Value = -770.3499999999999
$ws.Range("N110").Value = -8820.285599999999


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 40
$ws.Range("H40").Value = 42498.125
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 42498.125
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 42498.125
$ws.Range("N40").Value = -43028.125

# Row 107
$ws.Range("H107").Value = 10736.643
$ws.Range("I107").Value = 13615.857
$ws.Range("J107").Value = 7857.4287
$ws.Range("K107").Value = 13615.857
$ws.Range("L107").Value = 7857.4287
$ws.Range("M107").Value = -11695.857
$ws.Range("N107").Value = -11697.4287

# Row 134
$ws.Range("H134").Value = 972.2857
$ws.Range("I134").Value = 1023.5
$ws.Range("J134").Value = 665
$ws.Range("K134").Value = 3070.5
$ws.Range("L134").Value = 1995
$ws.Range("M134").Value = -535.5
$ws.Range("N134").Value = -7065


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1834.8334
$ws.Range("I58").Value = 1412.8572
$ws.Range("J58").Value = 2425.6
$ws.Range("K58").Value = 1412.8572
$ws.Range("L58").Value = 2425.6
$ws.Range("M58").Value = -1209.8572
$ws.Range("N58").Value = -2831.6

# Row 136
$ws.Range("H136").Value = 1834.8334
$ws.Range("I136").Value = 1412.8572
$ws.Range("J136").Value = 2425.6
$ws.Range("K136").Value = 4238.571599999999
$ws.Range("L136").Value = 7276.799999999999
$ws.Range("M136").Value = -1688.571599999999
$ws.Range("N136").Value = -12376.8


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 6
$ws.Range("H6").Value = 267.36365
$ws.Range("I6").Value = 341.16666
$ws.Range("J6").Value = 178.8
$ws.Range("K6").Value = 1023.49998
$ws.Range("L6").Value = 536.4000000000001
$ws.Range("M6").Value = -910.4999799999999
$ws.Range("N6").Value = -762.4000000000001

# Row 21
$ws.Range("H21").Value = 125.25
$ws.Range("I21").Value = 145.33333
$ws.Range("J21").Value = 65
$ws.Range("K21").Value = 435.99999
$ws.Range("L21").Value = 195
$ws.Range("M21").Value = -262.99999
$ws.Range("N21").Value = -541

# Row 122
$ws.Range("H122").Value = 23815726
$ws.Range("I122").Value = 733.3333
$ws.Range("J122").Value = 30310724
$ws.Range("K122").Value = 6599.9997
$ws.Range("L122").Value = 272796516
$ws.Range("M122").Value = -4149.9997
$ws.Range("N122").Value = -272801416


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 20800
$ws.Range("I43").Value = 20800
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 20800
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -20649
$ws.Range("N43").ClearContents()

# Row 97
$ws.Range("H97").Value = 13200.75
$ws.Range("I97").Value = 767.6667
$ws.Range("J97").Value = 50500
$ws.Range("K97").Value = 767.6667
$ws.Range("L97").Value = 50500
$ws.Range("M97").Value = -271.6667
$ws.Range("N97").Value = -51492


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 2901.5
$ws.Range("I132").Value = 2480.55
$ws.Range("J132").Value = 5006.25
$ws.Range("K132").Value = 7441.650000000001
$ws.Range("L132").Value = 15018.75
$ws.Range("M132").Value = -4911.650000000001
$ws.Range("N132").Value = -20078.75


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 109
$ws.Range("H109").Value = 17000
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 17000
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 17000
$ws.Range("N109").Value = -19774

# Row 115
$ws.Range("H115").Value = 28999.666
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 28999.666
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 28999.666
$ws.Range("N115").Value = -32133.666

# Row 126
$ws.Range("H126").Value = 4500.8945
$ws.Range("I126").Value = 6701.7
$ws.Range("J126").Value = 2055.5557
$ws.Range("K126").Value = 20105.1
$ws.Range("L126").Value = 6166.6671
$ws.Range("M126").Value = -17635.1
$ws.Range("N126").Value = -11106.6671

# Row 136
$ws.Range("H136").Value = 3810.4
$ws.Range("I136").Value = 3513.25
$ws.Range("J136").Value = 4999
$ws.Range("K136").Value = 10539.75
$ws.Range("L136").Value = 4999
$ws.Range("M136").Value = -7989.75
$ws.Range("N136").Value = -20097

